# Donor constraint analysis, refactored into a reusable function.
# Column J = "CET ($65.8) + Demand constraint + Drug budget + Donor constraint"
function Set-DonorConstraintResult($Worksheet, $Row, $Value) {
    $Worksheet.Cells.Item($Row, 10).Value = $Value
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimal_coverage")

$donorConstraintResults = @{
    2 = 0.4699999999999921
    7 = 0.161399429474903
    11 = 0.2689990491248384
    12 = 0.2420991442123545
    13 = 0.2420991442123545
    18 = 0.4000000000000001
    19 = 0.161399429474903
    20 = 0.161399429474903
    21 = 0.161399429474903
    22 = 0.06724976228120959
    23 = 0.06724976228120959
    24 = 0.06724976228120959
    26 = 0.2420991442123545
    30 = 0.1882993343873868
    31 = 0.1882993343873868
    38 = 0.2555490966685964
    39 = 0.2555490966685964
    40 = 0.2151992392998707
    41 = 0.2151992392998707
    42 = 0.2151992392998707
    43 = 0.2689990491248384
    45 = 0.9499999999999997
    50 = 0.161399429474903
    51 = 0.1882993343873868
    53 = 0.161399429474903
    54 = 0.161399429474903
    55 = 0.161399429474903
    58 = 0.5000000000003847
    65 = 0.9499999999993445
    66 = 0.9499999999994823
    69 = 0
    70 = 0
    95 = 0.1882993343873868
    96 = 0.1344995245624192
    101 = 0.2186871614598355
    107 = 0.2420991442123545
    108 = 0.6000000000002502
    109 = 0.2420991442123545
    122 = 0.2420991442123545
    123 = 0.2420991442123545
    127 = 0.2151992392998707
    128 = 0.2151992392998707
    129 = 0.2636190681423416
    130 = 0.2582390871598448
    131 = 0.2689990491248384
    132 = 0.2609290776510932
    134 = 0.2636190681423416
    135 = 0.2636190681423416
    138 = 0.161399429474903
    140 = 0.59
}

foreach ($row in $donorConstraintResults.Keys) {
    Set-DonorConstraintResult $ws $row $donorConstraintResults[$row]
}
